# Update the cryptos worksheet with the latest scraped coinranking.com data.
# (rows 2-51: price in column D, 1h volume/change in column E; a couple of rows
#  also had their coin name/link re-ranked into a new row position.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.345.82'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.577.63'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.47'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -3.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '23.87'
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.246'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0588'
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0897'
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').Value = '1.804.85'
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').Value = '1.570.67'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.68'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '28.405.06'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.516'
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.73'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '231.18'
$ws.Range('E19').Value = '  +1.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.43'
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.06'
$ws.Range('E24').Value = '  -1.04%  '
$ws.Range('E25').Value = '  +2.52%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.66'
$ws.Range('E26').Value = '  +0.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.02'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.104'
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.37'
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('E32').Value = '  -2.91%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('D35').Value = '1.390.86'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.09'
$ws.Range('E36').Value = '  +8.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.50'
$ws.Range('E37').Value = '  -2.77%  '
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.64'
$ws.Range('E39').Value = '  +3.40%  '
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.519'
$ws.Range('E41').Value = '  -2.43%  '
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('E43').Value = '  +2.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.784'
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('E45').Value = '  -3.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0456'
$ws.Range('E46').Value = '  -2.44%  '
$ws.Range('E47').Value = '  -4.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '62.51'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('D49').Value = '1.715.54'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '85.66'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '41.99'
$ws.Range('E51').Value = '  +5.91%  '
